$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 4623
$ws1.Range("F15").Value = 1028
$ws1.Range("F17").Value = 244
$ws1.Range("F22").Value = 3609
$ws1.Range("F23").Value = 5947
$ws1.Range("F45").Value = 458

# Sheet "全部类型" (All Types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 4623
$ws4.Range("F16").Value = 1028
$ws4.Range("F18").Value = 244
$ws4.Range("F23").Value = 3609
$ws4.Range("F24").Value = 5947
$ws4.Range("F46").Value = 458
